$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced to
# text format first, otherwise Excel auto-converts them to numbers and the
# exact source formatting (e.g. trailing zeros) is lost.
$numericLikeRange = $ws.Range("D4,D5,D8,D9,D10,D11,D12,D15,D17,D20,D21,D22,D23,D25,D27,D28,D29,D30,D31,D32,D33,D35,D36,D37,D38,D40,D42,D43,D44,D47,D48,D49,D50,D51")
foreach ($cell in $numericLikeRange) {
    $cell.NumberFormat = "@"
}

$ws.Range("D2").Value = "25.816.86"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.635.60"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "215.28"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "0.2571"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "0.06411"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").Value = "20.02"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("D11").Value = "0.07782"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "4.287"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "1.862.29"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "1.634.43"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "0.5610"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "0.0₅7632"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "63.03"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").Value = "25.851.08"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "194.12"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "4.329"
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("D22").Value = "9.881"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "6.096"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D25").Value = "1.779"
$ws.Range("E25").Value = "  -6.71%  "
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").Value = "0.1256"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").Value = "6.806"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "15.42"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").Value = "1.241"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "0.04906"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").Value = "3.303"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").Value = "3.233"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("D35").Value = "2.380"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "0.9027"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "2.576"
$ws.Range("D38").Value = "0.5522"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "1.126.65"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").Value = "0.01557"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D42").Value = "5.516"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "0.7998"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "98.06"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "1.772.48"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("E46").Value = "  -7.35%  "
$ws.Range("D47").Value = "55.42"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("D48").Value = "0.4262"
$ws.Range("E48").Value = "  -4.38%  "
$ws.Range("D49").Value = "7.741"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("D50").Value = "0.05031"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").Value = "0.9998"
$ws.Range("E51").Value = "  +0.48%  "
